# Update the worker ranking-matrix simulation results.
# The underlying "matrices" random draw was re-run, which changes the
# mat_rank score (column G) for every worker, reshuffles a handful of the
# "matrices" counts (column C) and, because three pairs of workers ended up
# swapping places in the (re-simulated) ranking, their identity columns
# (prolificid / name / race) move together as a unit for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- matrices = female group (rows 2-13) ---
$ws.Range("G2").Value = 13.25581603006527
$ws.Range("G3").Value = 13.03564410204013
$ws.Range("G4").Value = 8.387878449008936
$ws.Range("G5").Value = 8.356292063322577
$ws.Range("G6").Value = 8.180015286402934
$ws.Range("G7").Value = 5.393336665672788

# Row 8 and row 9 swap identities (prolificid/name/race) plus new scores
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("E8").Value = "Shadaisia"
$ws.Range("G8").Value = 5.339669197139461
$ws.Range("H8").Value = "Black or African American"

$ws.Range("C9").Value = 32
$ws.Range("D9").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("E9").Value = "Kellie"
$ws.Range("G9").Value = 5.108019693417147
$ws.Range("H9").Value = "White"

$ws.Range("G10").Value = 4.268640122598316
$ws.Range("G11").Value = 4.254495598246366
$ws.Range("G12").Value = 2.420025270519735
$ws.Range("G13").Value = 1.496024677253027

# --- matrices = male group (rows 14-25) ---
$ws.Range("G14").Value = 14.35604799398173
$ws.Range("G15").Value = 13.17756464437572
$ws.Range("G16").Value = 8.158367614863963
$ws.Range("G17").Value = 7.132419507397405
$ws.Range("G18").Value = 6.306267974076017
$ws.Range("G19").Value = 6.305467982787811

# Rows 20, 21 and 22 rotate identities (prolificid/name/race) plus new scores
$ws.Range("C20").Value = 30
$ws.Range("D20").Value = "60c2341fe95d71ee52c043f0"
$ws.Range("E20").Value = "Matthew"
$ws.Range("G20").Value = 5.443833869706829
$ws.Range("H20").Value = "White"

$ws.Range("C21").Value = 32
$ws.Range("D21").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("E21").Value = "Jamarii"
$ws.Range("G21").Value = 5.210446373867417
$ws.Range("H21").Value = "Black or African American"

$ws.Range("C22").Value = 33
$ws.Range("D22").Value = "60b322994d0b901954690036"
$ws.Range("E22").Value = "Brennan"
$ws.Range("G22").Value = 5.049269166493271

$ws.Range("G23").Value = 3.295791998891051
$ws.Range("G24").Value = 1.043195326962711
$ws.Range("G25").Value = 0.2516358054655306
